# "Start on the outside-in presentation."
#
# 1) The datetimeFigureOut auto-date fields on the slide master, all
#    (non-title) slide layouts, the handout master and the notes master
#    rolled from 11/2/2010 to 11/3/2010 (opened the deck a day later).
# 2) The title-slide subtitle's literal date line was bumped a month,
#    from 10/6/2010 to 11/6/2010.
# 3) Slide 13's two adjacent runs ("It's hard to go wrong ... use it" + ".")
#    collapsed back into a single run once the trailing period was
#    retyped as part of the same text.

$p = $ppt.ActivePresentation

# --- 1a) Slide master date placeholder ---------------------------------
$master = $p.SlideMaster
$master.Shapes.Item(4).TextFrame.TextRange.Text = "11/3/2010"

# --- 1b) Handout master date placeholder --------------------------------
$handout = $p.HandoutMaster
$handout.Shapes.Item(2).TextFrame.TextRange.Text = "11/3/2010"

# --- 1c) Notes master date placeholder -----------------------------------
$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item(2).TextFrame.TextRange.Text = "11/3/2010"

# --- 1d) Every slide layout's date placeholder (Title Slide has none) ---
$layouts = $master.CustomLayouts
$dateShapeByLayout = @{2=4; 3=3; 4=4; 5=6; 6=4; 7=4; 8=3; 9=3; 10=2; 11=1; 12=2}
foreach ($layoutIdx in $dateShapeByLayout.Keys) {
    $layout = $layouts.Item($layoutIdx)
    $shapeIdx = $dateShapeByLayout[$layoutIdx]
    $layout.Shapes.Item($shapeIdx).TextFrame.TextRange.Text = "11/3/2010"
}

# --- 2) Title slide subtitle: "10/6/2010" -> "11/6/2010" ----------------
$titleSlide = $p.Slides.Item(1)
$subtitle = $titleSlide.Shapes.Item(2)
$subtitleRange = $subtitle.TextFrame.TextRange
$dateText = $subtitleRange.Characters(16, 9)
$dateText.Text = "11/6/2010"

# --- 3) Slide 13: retype the sentence so the trailing "." merges back ---
#        into the first run instead of staying a separate run.
$slide13 = $p.Slides.Item(13)
$content = $slide13.Shapes.Item(2)
$contentRange = $content.TextFrame.TextRange
$sentence = $contentRange.Characters(70, 69)
$sentence.Text = "It’s hard to go wrong with TDD, but when BDD is a better fit, use it."
